$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): add new headers in columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style of an existing header cell (AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill data rows 2-48 with the team's win/loss/tie record
$lastRow = 48
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 85   # AD = column 30 -> Wins
    $ws.Cells.Item($r, 31).Value = 77   # AE = column 31 -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF = column 32 -> Ties
}
